$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the status column (C) for the first three data rows:
# "da xu ly" -> "da mo lop". Leading apostrophe keeps the quote-prefix
# cell style (s="3") intact instead of letting the write reset it.
$ws.Range("C2").Value = "'da mo lop"
$ws.Range("C3").Value = "'da mo lop"
$ws.Range("C4").Value = "'da mo lop"

# Move the active selection from A5 to D2.
$ws.Range("D2").Select()

# Persist so the shared-string table is rebuilt/compacted (drops the
# now-unused "da xu ly" entry), matching canonical Excel save behavior.
$wb.Save()
